# Connected FTT:H2 with CLEAFS and included a split market.
# Adds new "WDM" mandate variables (split-market demand mandates) plus a
# "HYGR" (endogenous green-fertiliser-driven hydrogen demand growth) row to
# the FTT-H2 variable-definitions sheet, registers their time horizons on
# the Time_Horizons sheet, and leaves the workbook focused on FTT-H2.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. FTT-H2 sheet: append six new variable-definition rows (26-31)
# ---------------------------------------------------------------------
$ftth2 = $wb.Worksheets.Item("FTT-H2")

# Row 26: WDM1 - Mandate for NH3 for fertiliser
$ftth2.Range("A26").Value = "WDM1"
$ftth2.Range("B26").Value = 1
$ftth2.Range("C26").Formula = "=C25+1000"
$ftth2.Range("D26").Value = "Mandate for NH3 for fertiliser"
$ftth2.Range("E26").Value = "RSHORTTI"
$ftth2.Range("F26").Value = "TIME"
$ftth2.Range("G26").Value = 0
$ftth2.Range("H26").Value = 0
$ftth2.Range("I26").Value = "S1"

# Row 27: WDM2 - Mandate for NH3 for chemicals
$ftth2.Range("A27").Value = "WDM2"
$ftth2.Range("B27").Value = 1
$ftth2.Range("C27").Formula = "=C26+1000"
$ftth2.Range("D27").Value = "Mandate for NH3 for chemicals"
$ftth2.Range("E27").Value = "RSHORTTI"
$ftth2.Range("F27").Value = "TIME"
$ftth2.Range("G27").Value = 0
$ftth2.Range("H27").Value = 0
$ftth2.Range("I27").Value = "S2"

# Row 28: WDM3 - Mandate for MeOH for chemicals
$ftth2.Range("A28").Value = "WDM3"
$ftth2.Range("B28").Value = 1
$ftth2.Range("C28").Formula = "=C27+1000"
$ftth2.Range("D28").Value = "Mandate for MeOH for chemicals"
$ftth2.Range("E28").Value = "RSHORTTI"
$ftth2.Range("F28").Value = "TIME"
$ftth2.Range("G28").Value = 0
$ftth2.Range("H28").Value = 0
$ftth2.Range("I28").Value = "S3"

# Row 29: WDM4 - Mandate for H2 for oil refining
$ftth2.Range("A29").Value = "WDM4"
$ftth2.Range("B29").Value = 1
$ftth2.Range("C29").Formula = "=C28+1000"
$ftth2.Range("D29").Value = "Mandate for H2 for oil refining"
$ftth2.Range("E29").Value = "RSHORTTI"
$ftth2.Range("F29").Value = "TIME"
$ftth2.Range("G29").Value = 0
$ftth2.Range("H29").Value = 0
$ftth2.Range("I29").Value = "S4"

# Row 30: WDM5 - Mandate for H2 for energy purposes
$ftth2.Range("A30").Value = "WDM5"
$ftth2.Range("B30").Value = 1
$ftth2.Range("C30").Formula = "=C29+1000"
$ftth2.Range("D30").Value = "Mandate for H2 for energy purposes"
$ftth2.Range("E30").Value = "RSHORTTI"
$ftth2.Range("F30").Value = "TIME"
$ftth2.Range("G30").Value = 0
$ftth2.Range("H30").Value = 0
$ftth2.Range("I30").Value = "S5"

# Row 31: HYGR - endogenous green-fertiliser-driven hydrogen demand growth
$ftth2.Range("A31").Value = "HYGR"
$ftth2.Range("B31").Value = 1
$ftth2.Range("C31").Value = 3609000
$ftth2.Range("D31").Value = "FTT:Hydrogen global cumulative capacity (kt H2)"
$ftth2.Range("E31").Value = "HYTI"
$ftth2.Range("F31").Value = 0
$ftth2.Range("G31").Value = 0
$ftth2.Range("H31").Value = "TIME"
$ftth2.Range("I31").Value = "S0"

# ---------------------------------------------------------------------
# 2. Time_Horizons sheet: update PFRA horizon, add horizons for new vars
# ---------------------------------------------------------------------
$th = $wb.Worksheets.Item("Time_Horizons")

# PFRA's horizon moves from tl_2010 to tl_2020 (now tracked up to 2020)
$th.Range("B68").Value = "tl_2020"

$th.Range("A91").Value = "WDM1"
$th.Range("B91").Value = "tl_2001"

$th.Range("A92").Value = "WDM2"
$th.Range("B92").Value = "tl_2001"

$th.Range("A93").Value = "WDM3"
$th.Range("B93").Value = "tl_2001"

$th.Range("A94").Value = "WDM4"
$th.Range("B94").Value = "tl_2001"

$th.Range("A95").Value = "WDM5"
$th.Range("B95").Value = "tl_2001"

$th.Range("A96").Value = "HYGR"
$th.Range("B96").Value = "tl_2022"

# ---------------------------------------------------------------------
# 3. View state: update Time_Horizons' scroll/selection, then leave the
#    workbook focused on the FTT-H2 sheet (mirrors the author's session)
# ---------------------------------------------------------------------
$th.Activate()
$th.Range("E97").Select()
$excel.ActiveWindow.ScrollRow = 75

$ftth2.Activate()
$ftth2.Range("A6").Select()
